# Fix the source link on the "Źródła" slide (last slide, #12):
# change https://blockly-games.appspot.com/bird -> .../movie
# and merge away the trailing empty paragraph, per the commit's intent
# ("Poprawiono link do źródła").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Setting TextRange.Text directly tries to keep unchanged runs/characters
# intact by diffing against the previous text, which would split the
# hyperlink run in two (shared "https://blockly-games.appspot.com/" prefix
# kept, "bird"/"movie" suffix turned into a new run). First assign an
# unrelated placeholder string so the following assignment has nothing in
# common with it, forcing the whole paragraph to be rewritten as a single
# run (preserving the single <a:r> with its hlinkClick, plus endParaRPr).
$tr.Text = "placeholder-no-overlap-xxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxxx`r"
$tr.Text = "https://blockly-games.appspot.com/movie`r"

# The trailing `r above keeps the second (empty) paragraph around; remove
# it now so only the single, corrected paragraph remains.
$para2 = $tr.Paragraphs(2, 1)
$para2.Delete()
